$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates (K1, L1) ---
$ws.Range('K1').Value = 'Type (for ADD)'
$ws.Range('L1').Value = 'Email/Domain (for ADD)'

# --- Data rows 2-11 (project structure reshuffled + refreshed) ---
# Row 2
$ws.Range('A2').Value = '/Python_Admin_Tool_TESTING/Sub-Folder 1'
$ws.Range('B2').Value = 'Sub-Folder 1'
$ws.Range('C2').Value = '1jPqgww8lNGleK7h15iHuNdyUbZDf0idz'
$ws.Range('D2').Value = 'Commenter'
$ws.Range('E2').Value = 'domain'
$ws.Range('F2').Value = 'bioaccessla.com'
$ws.Range('G2').Value = 'jmoreno@bioaccessla.com'
$ws.Range('H2').Value = 'https://drive.google.com/drive/folders/1jPqgww8lNGleK7h15iHuNdyUbZDf0idz'

# Row 3
$ws.Range('A3').Value = '/Python_Admin_Tool_TESTING/Sub-Folder 1'
$ws.Range('B3').Value = 'Sub-Folder 1'
$ws.Range('C3').Value = '1jPqgww8lNGleK7h15iHuNdyUbZDf0idz'
$ws.Range('D3').Value = 'Viewer'
$ws.Range('E3').Value = 'user'
$ws.Range('F3').Value = 'ernie.moreno62@gmail.com'
$ws.Range('G3').Value = 'jmoreno@bioaccessla.com'
$ws.Range('H3').Value = 'https://drive.google.com/drive/folders/1jPqgww8lNGleK7h15iHuNdyUbZDf0idz'

# Row 4
$ws.Range('A4').Value = '/Python_Admin_Tool_TESTING/Sub-Folder 1'
$ws.Range('B4').Value = 'Sub-Folder 1'
$ws.Range('C4').Value = '1jPqgww8lNGleK7h15iHuNdyUbZDf0idz'
$ws.Range('D4').Value = 'Owner'
$ws.Range('E4').Value = 'user'
$ws.Range('F4').Value = 'jmoreno@bioaccessla.com'
$ws.Range('G4').Value = 'jmoreno@bioaccessla.com'
$ws.Range('H4').Value = 'https://drive.google.com/drive/folders/1jPqgww8lNGleK7h15iHuNdyUbZDf0idz'

# Row 5
$ws.Range('A5').Value = '/Python_Admin_Tool_TESTING/Sub-Folder 1/Test Sheet 2'
$ws.Range('B5').Value = 'Test Sheet 2'
$ws.Range('C5').Value = '1Wan1C_Cxndc2M6yXKa8vxJLXCgYoBXFzDVwRYAfGbVY'
$ws.Range('D5').Value = 'Editor'
$ws.Range('E5').Value = 'group'
$ws.Range('F5').Value = 'jesus_test_group@bioaccessla.com'
$ws.Range('G5').Value = 'jmoreno@bioaccessla.com'
$ws.Range('H5').Value = 'https://docs.google.com/spreadsheets/d/1Wan1C_Cxndc2M6yXKa8vxJLXCgYoBXFzDVwRYAfGbVY/edit?usp=drivesdk'

# Row 6
$ws.Range('A6').Value = '/Python_Admin_Tool_TESTING/Sub-Folder 1/Test Sheet 2'
$ws.Range('B6').Value = 'Test Sheet 2'
$ws.Range('C6').Value = '1Wan1C_Cxndc2M6yXKa8vxJLXCgYoBXFzDVwRYAfGbVY'
$ws.Range('D6').Value = 'Commenter'
$ws.Range('E6').Value = 'domain'
$ws.Range('F6').Value = 'bioaccessla.com'
$ws.Range('G6').Value = 'jmoreno@bioaccessla.com'
$ws.Range('H6').Value = 'https://docs.google.com/spreadsheets/d/1Wan1C_Cxndc2M6yXKa8vxJLXCgYoBXFzDVwRYAfGbVY/edit?usp=drivesdk'

# Row 7
$ws.Range('A7').Value = '/Python_Admin_Tool_TESTING/Sub-Folder 1/Test Sheet 2'
$ws.Range('B7').Value = 'Test Sheet 2'
$ws.Range('C7').Value = '1Wan1C_Cxndc2M6yXKa8vxJLXCgYoBXFzDVwRYAfGbVY'
$ws.Range('D7').Value = 'Viewer'
$ws.Range('E7').Value = 'user'
$ws.Range('F7').Value = 'ernie.moreno62@gmail.com'
$ws.Range('G7').Value = 'jmoreno@bioaccessla.com'
$ws.Range('H7').Value = 'https://docs.google.com/spreadsheets/d/1Wan1C_Cxndc2M6yXKa8vxJLXCgYoBXFzDVwRYAfGbVY/edit?usp=drivesdk'

# Row 8
$ws.Range('A8').Value = '/Python_Admin_Tool_TESTING/Sub-Folder 1/Test Sheet 2'
$ws.Range('B8').Value = 'Test Sheet 2'
$ws.Range('C8').Value = '1Wan1C_Cxndc2M6yXKa8vxJLXCgYoBXFzDVwRYAfGbVY'
$ws.Range('D8').Value = 'Owner'
$ws.Range('E8').Value = 'user'
$ws.Range('F8').Value = 'jmoreno@bioaccessla.com'
$ws.Range('G8').Value = 'jmoreno@bioaccessla.com'
$ws.Range('H8').Value = 'https://docs.google.com/spreadsheets/d/1Wan1C_Cxndc2M6yXKa8vxJLXCgYoBXFzDVwRYAfGbVY/edit?usp=drivesdk'

# Row 9
$ws.Range('A9').Value = '/Python_Admin_Tool_TESTING/Test Doc 1'
$ws.Range('B9').Value = 'Test Doc 1'
$ws.Range('C9').Value = '1O90b5jSuK3lIz-RYZIEtAAlA3-IQ_vmxgulyB_6vY2U'
$ws.Range('D9').Value = 'Viewer'
$ws.Range('E9').Value = 'domain'
$ws.Range('F9').Value = 'bioaccessla.com'
$ws.Range('G9').Value = 'jmoreno@bioaccessla.com'
$ws.Range('H9').Value = 'https://docs.google.com/document/d/1O90b5jSuK3lIz-RYZIEtAAlA3-IQ_vmxgulyB_6vY2U/edit?usp=drivesdk'

# Row 10
$ws.Range('A10').Value = '/Python_Admin_Tool_TESTING/Test Doc 1'
$ws.Range('B10').Value = 'Test Doc 1'
$ws.Range('C10').Value = '1O90b5jSuK3lIz-RYZIEtAAlA3-IQ_vmxgulyB_6vY2U'
$ws.Range('D10').Value = 'Commenter'
$ws.Range('E10').Value = 'user'
$ws.Range('F10').Value = 'ernie.moreno62@gmail.com'
$ws.Range('G10').Value = 'jmoreno@bioaccessla.com'
$ws.Range('H10').Value = 'https://docs.google.com/document/d/1O90b5jSuK3lIz-RYZIEtAAlA3-IQ_vmxgulyB_6vY2U/edit?usp=drivesdk'

# Row 11
$ws.Range('A11').Value = '/Python_Admin_Tool_TESTING/Test Doc 1'
$ws.Range('B11').Value = 'Test Doc 1'
$ws.Range('C11').Value = '1O90b5jSuK3lIz-RYZIEtAAlA3-IQ_vmxgulyB_6vY2U'
$ws.Range('D11').Value = 'Owner'
$ws.Range('E11').Value = 'user'
$ws.Range('F11').Value = 'jmoreno@bioaccessla.com'
$ws.Range('G11').Value = 'jmoreno@bioaccessla.com'
$ws.Range('H11').Value = 'https://docs.google.com/document/d/1O90b5jSuK3lIz-RYZIEtAAlA3-IQ_vmxgulyB_6vY2U/edit?usp=drivesdk'

# --- Conditional formatting: highlight rows by Action_Type (ADD/REMOVE/MODIFY) ---
# Adds 3 dxf fills (solid green/red/yellow) + 3 expression cfRules on A2:K1048576.
# xlExpression = 2, xlEqual (unused for expression rules, required positional arg) = 3
$cfRange = $ws.Range('A2:K1048576')

# ADD -> light green FFD8E9BB (RGB 216,233,187 => 187*65536 + 233*256 + 216)
$cfAdd = $cfRange.FormatConditions.Add(2, 3, '=$I2="ADD"')
$cfAdd.Interior.Color = 12315096

# REMOVE -> light red FFFFC7CE (RGB 255,199,206 => 206*65536 + 199*256 + 255)
$cfRemove = $cfRange.FormatConditions.Add(2, 3, '=$I2="REMOVE"')
$cfRemove.Interior.Color = 13551615

# MODIFY -> light yellow FFFFEB9C (RGB 255,235,156 => 156*65536 + 235*256 + 255)
$cfModify = $cfRange.FormatConditions.Add(2, 3, '=$I2="MODIFY"')
$cfModify.Interior.Color = 10284031
